$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: add P1/Q1 headers (14, 15), matching the bold/bordered style
#     used by the rest of row 1 (copy O1's formatting, then overwrite value) ---
$ws.Range("O1").Copy($ws.Range("P1"))
$ws.Range("O1").Copy($ws.Range("Q1"))
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# --- Rows 2-25: swap I<->K and M<->O, then append new P/Q columns (value 2) ---
for ($r = 2; $r -le 25; $r++) {
    $i = $ws.Cells.Item($r, 9).Value()
    $k = $ws.Cells.Item($r, 11).Value()
    $m = $ws.Cells.Item($r, 13).Value()
    $o = $ws.Cells.Item($r, 15).Value()

    $ws.Cells.Item($r, 9).Value = $k
    $ws.Cells.Item($r, 11).Value = $i
    $ws.Cells.Item($r, 13).Value = $o
    $ws.Cells.Item($r, 15).Value = $m

    $ws.Cells.Item($r, 16).Value = 2
    $ws.Cells.Item($r, 17).Value = 2
}
